$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 12.55295333333333
$ws.Range("H2").Value = 37.65886
$ws.Range("I2").Value = 0.3363704472878066
$ws.Range("J2").Value = 0.3591006154861918
$ws.Range("M2").Value = 38.10639333333333
$ws.Range("N2").Value = 114.31918
$ws.Range("O2").Value = 0.3831479157160237
$ws.Range("P2").Value = 0.4159903984418967
$ws.Range("Q2").Value = 478.3477772149778
$ws.Range("R2").Value = 4305.1299949348
$ws.Range("S2").Value = 0.1288796357867897
$ws.Range("T2").Value = 0.1493824081168313

# Row 3
$ws.Range("G3").Value = 12.55295333333333
$ws.Range("H3").Value = 37.65886
$ws.Range("I3").Value = 0.3363704472878066
$ws.Range("J3").Value = 0.3591006154861918
$ws.Range("O3").Value = 0.09199521176963764
$ws.Range("P3").Value = 0.09988081163714851
$ws.Range("Q3").Value = 114.85304567608
$ws.Range("R3").Value = 1033.67741108472
$ws.Range("S3").Value = 0.0309444705312895
$ws.Range("T3").Value = 0.03586726093416042

# Row 4
$ws.Range("G4").Value = 12.55295333333333
$ws.Range("H4").Value = 37.65886
$ws.Range("I4").Value = 0.3363704472878066
$ws.Range("J4").Value = 0.3591006154861918
$ws.Range("M4").Value = 15.023598
$ws.Range("N4").Value = 45.070794
$ws.Range("O4").Value = 0.1510575983904562
$ws.Range("P4").Value = 0.1640058785774412
$ws.Range("Q4").Value = 188.59052459276
$ws.Range("R4").Value = 1697.31472133484
$ws.Range("S4").Value = 0.05081131193681961
$ws.Range("T4").Value = 0.05889461194051279

# Row 5
$ws.Range("G5").Value = 12.55295333333333
$ws.Range("H5").Value = 37.65886
$ws.Range("I5").Value = 0.3363704472878066
$ws.Range("J5").Value = 0.3591006154861918
$ws.Range("M5").Value = 23.556204
$ws.Range("N5").Value = 47.112408
$ws.Range("O5").Value = 0.236850294013169
$ws.Range("P5").Value = 0.1714350065796238
$ws.Range("Q5").Value = 295.69992952248
$ws.Range("R5").Value = 1774.19957713488
$ws.Range("S5").Value = 0.07966943933745815
$ws.Range("T5").Value = 0.06156241637862227

# Row 6
$ws.Range("G6").Value = 12.55295333333333
$ws.Range("H6").Value = 37.65886
$ws.Range("I6").Value = 0.3363704472878066
$ws.Range("J6").Value = 0.3591006154861918
$ws.Range("M6").Value = 13.62041
$ws.Range("N6").Value = 40.86123000000001
$ws.Range("O6").Value = 0.1369489801107134
$ws.Range("P6").Value = 0.1486879047638899
$ws.Range("Q6").Value = 170.9763711108667
$ws.Range("R6").Value = 1538.7873399978
$ws.Range("S6").Value = 0.04606558969544959
$ws.Range("T6").Value = 0.05339391811606513

# Row 7
$ws.Range("H7").Value = 39.447015
$ws.Range("I7").Value = 0.3523423194360853
$ws.Range("J7").Value = 0.3761517838190811
$ws.Range("M7").Value = 38.10639333333333
$ws.Range("N7").Value = 114.31918
$ws.Range("O7").Value = 0.3831479157160237
$ws.Range("P7").Value = 0.4159903984418967
$ws.Range("Q7").Value = 501.0611564719666
$ws.Range("R7").Value = 4509.5504082477
$ws.Range("S7").Value = 0.1349992253104856
$ws.Range("T7").Value = 0.1564755304255297

# Row 8
$ws.Range("H8").Value = 39.447015
$ws.Range("I8").Value = 0.3523423194360853
$ws.Range("J8").Value = 0.3761517838190811
$ws.Range("O8").Value = 0.09199521176963764
$ws.Range("P8").Value = 0.09988081163714851
$ws.Range("S8").Value = 0.03241380629192799
$ws.Range("T8").Value = 0.03757034546661105

# Row 9
$ws.Range("H9").Value = 39.447015
$ws.Range("I9").Value = 0.3523423194360853
$ws.Range("J9").Value = 0.3761517838190811
$ws.Range("M9").Value = 15.023598
$ws.Range("N9").Value = 45.070794
$ws.Range("O9").Value = 0.1510575983904562
$ws.Range("P9").Value = 0.1640058785774412
$ws.Range("Q9").Value = 197.54536521999
$ws.Range("R9").Value = 1777.90828697991
$ws.Range("S9").Value = 0.053223984585338
$ws.Range("T9").Value = 0.06169110378372013

# Row 10
$ws.Range("H10").Value = 39.447015
$ws.Range("I10").Value = 0.3523423194360853
$ws.Range("J10").Value = 0.3761517838190811
$ws.Range("M10").Value = 23.556204
$ws.Range("N10").Value = 47.112408
$ws.Range("O10").Value = 0.236850294013169
$ws.Range("P10").Value = 0.1714350065796238
$ws.Range("Q10").Value = 309.74064417702
$ws.Range("R10").Value = 1858.44386506212
$ws.Range("S10").Value = 0.08345238195171872
$ws.Range("T10").Value = 0.06448558353396142

# Row 11
$ws.Range("H11").Value = 39.447015
$ws.Range("I11").Value = 0.3523423194360853
$ws.Range("J11").Value = 0.3761517838190811
$ws.Range("M11").Value = 13.62041
$ws.Range("N11").Value = 40.86123000000001
$ws.Range("O11").Value = 0.1369489801107134
$ws.Range("P11").Value = 0.1486879047638899
$ws.Range("Q11").Value = 179.09483919205
$ws.Range("R11").Value = 1611.85355272845
$ws.Range("S11").Value = 0.04825292129661507
$ws.Range("T11").Value = 0.05592922060925883

# Row 12
$ws.Range("G12").Value = 2.132104
$ws.Range("H12").Value = 6.396312
$ws.Range("I12").Value = 0.05713211521624299
$ws.Range("J12").Value = 0.06099280689967021
$ws.Range("M12").Value = 38.10639333333333
$ws.Range("N12").Value = 114.31918
$ws.Range("O12").Value = 0.3831479157160237
$ws.Range("P12").Value = 0.4159903984418967
$ws.Range("Q12").Value = 81.24679365157333
$ws.Range("R12").Value = 731.2211428641599
$ws.Range("S12").Value = 0.02189005086555123
$ws.Range("T12").Value = 0.02537242204428347

# Row 13
$ws.Range("G13").Value = 2.132104
$ws.Range("H13").Value = 6.396312
$ws.Range("I13").Value = 0.05713211521624299
$ws.Range("J13").Value = 0.06099280689967021
$ws.Range("O13").Value = 0.09199521176963764
$ws.Range("P13").Value = 0.09988081163714851
$ws.Range("Q13").Value = 19.507651434336
$ws.Range("R13").Value = 175.568862909024
$ws.Range("S13").Value = 0.005255881038165611
$ws.Range("T13").Value = 0.006092011057166932

# Row 14
$ws.Range("G14").Value = 2.132104
$ws.Range("H14").Value = 6.396312
$ws.Range("I14").Value = 0.05713211521624299
$ws.Range("J14").Value = 0.06099280689967021
$ws.Range("M14").Value = 15.023598
$ws.Range("N14").Value = 45.070794
$ws.Range("O14").Value = 0.1510575983904562
$ws.Range("P14").Value = 0.1640058785774412
$ws.Range("Q14").Value = 32.031873390192
$ws.Range("R14").Value = 288.286860511728
$ws.Range("S14").Value = 0.008630240115532505
$ws.Range("T14").Value = 0.01000317888248463

# Row 15
$ws.Range("G15").Value = 2.132104
$ws.Range("H15").Value = 6.396312
$ws.Range("I15").Value = 0.05713211521624299
$ws.Range("J15").Value = 0.06099280689967021
$ws.Range("M15").Value = 23.556204
$ws.Range("N15").Value = 47.112408
$ws.Range("O15").Value = 0.236850294013169
$ws.Range("P15").Value = 0.1714350065796238
$ws.Range("Q15").Value = 50.224276773216
$ws.Range("R15").Value = 301.345660639296
$ws.Range("S15").Value = 0.0135317582865614
$ws.Range("T15").Value = 0.01045630225215469

# Row 16
$ws.Range("G16").Value = 2.132104
$ws.Range("H16").Value = 6.396312
$ws.Range("I16").Value = 0.05713211521624299
$ws.Range("J16").Value = 0.06099280689967021
$ws.Range("M16").Value = 13.62041
$ws.Range("N16").Value = 40.86123000000001
$ws.Range("O16").Value = 0.1369489801107134
$ws.Range("P16").Value = 0.1486879047638899
$ws.Range("Q16").Value = 29.04013064264
$ws.Range("R16").Value = 261.36117578376
$ws.Range("S16").Value = 0.007824184910432247
$ws.Range("T16").Value = 0.00906889266358049

# Row 17
$ws.Range("G17").Value = 7.086566
$ws.Range("H17").Value = 14.173132
$ws.Range("I17").Value = 0.1898924748509033
$ws.Range("J17").Value = 0.135149614846733
$ws.Range("M17").Value = 38.10639333333333
$ws.Range("N17").Value = 114.31918
$ws.Range("O17").Value = 0.3831479157160237
$ws.Range("P17").Value = 0.4159903984418967
$ws.Range("Q17").Value = 270.0434713786266
$ws.Range("R17").Value = 1620.26082827176
$ws.Range("S17").Value = 0.07275690594928104
$ws.Range("T17").Value = 0.05622094212936135

# Row 18
$ws.Range("G18").Value = 7.086566
$ws.Range("H18").Value = 14.173132
$ws.Range("I18").Value = 0.1898924748509033
$ws.Range("J18").Value = 0.135149614846733
$ws.Range("O18").Value = 0.09199521176963764
$ws.Range("P18").Value = 0.09988081163714851
$ws.Range("Q18").Value = 64.838422231944
$ws.Range("R18").Value = 389.030533391664
$ws.Range("S18").Value = 0.01746919843736943
$ws.Range("T18").Value = 0.01349885322333971

# Row 19
$ws.Range("G19").Value = 7.086566
$ws.Range("H19").Value = 14.173132
$ws.Range("I19").Value = 0.1898924748509033
$ws.Range("J19").Value = 0.135149614846733
$ws.Range("M19").Value = 15.023598
$ws.Range("N19").Value = 45.070794
$ws.Range("O19").Value = 0.1510575983904562
$ws.Range("P19").Value = 0.1640058785774412
$ws.Range("Q19").Value = 106.465718784468
$ws.Range("R19").Value = 638.794312706808
$ws.Range("S19").Value = 0.02868470120339755
$ws.Range("T19").Value = 0.02216533132234125

# Row 20
$ws.Range("G20").Value = 7.086566
$ws.Range("H20").Value = 14.173132
$ws.Range("I20").Value = 0.1898924748509033
$ws.Range("J20").Value = 0.135149614846733
$ws.Range("M20").Value = 23.556204
$ws.Range("N20").Value = 47.112408
$ws.Range("O20").Value = 0.236850294013169
$ws.Range("P20").Value = 0.1714350065796238
$ws.Range("Q20").Value = 166.932594355464
$ws.Range("R20").Value = 667.730377421856
$ws.Range("S20").Value = 0.04497608849932474
$ws.Range("T20").Value = 0.02316937511048331

# Row 21
$ws.Range("G21").Value = 7.086566
$ws.Range("H21").Value = 14.173132
$ws.Range("I21").Value = 0.1898924748509033
$ws.Range("J21").Value = 0.135149614846733
$ws.Range("M21").Value = 13.62041
$ws.Range("N21").Value = 40.86123000000001
$ws.Range("O21").Value = 0.1369489801107134
$ws.Range("P21").Value = 0.1486879047638899
$ws.Range("Q21").Value = 96.52193441206002
$ws.Range("R21").Value = 579.1316064723601
$ws.Range("S21").Value = 0.02600558076153049
$ws.Range("T21").Value = 0.02009511306120744

# Row 22
$ws.Range("G22").Value = 2.398207
$ws.Range("H22").Value = 7.194621
$ws.Range("I22").Value = 0.06426264320896187
$ws.Range("J22").Value = 0.06860517894832399
$ws.Range("M22").Value = 38.10639333333333
$ws.Range("N22").Value = 114.31918
$ws.Range("O22").Value = 0.3831479157160237
$ws.Range("P22").Value = 0.4159903984418967
$ws.Range("Q22").Value = 91.38701923675332
$ws.Range("R22").Value = 822.4831731307798
$ws.Range("S22").Value = 0.02462209780391623
$ws.Range("T22").Value = 0.02853909572589092

# Row 23
$ws.Range("G23").Value = 2.398207
$ws.Range("H23").Value = 7.194621
$ws.Range("I23").Value = 0.06426264320896187
$ws.Range("J23").Value = 0.06860517894832399
$ws.Range("O23").Value = 0.09199521176963764
$ws.Range("P23").Value = 0.09988081163714851
$ws.Range("Q23").Value = 21.942356575188
$ws.Range("R23").Value = 197.481209176692
$ws.Range("S23").Value = 0.005911855470885114
$ws.Range("T23").Value = 0.006852340955870415

# Row 24
$ws.Range("G24").Value = 2.398207
$ws.Range("H24").Value = 7.194621
$ws.Range("I24").Value = 0.06426264320896187
$ws.Range("J24").Value = 0.06860517894832399
$ws.Range("M24").Value = 15.023598
$ws.Range("N24").Value = 45.070794
$ws.Range("O24").Value = 0.1510575983904562
$ws.Range("P24").Value = 0.1640058785774412
$ws.Range("Q24").Value = 36.029697888786
$ws.Range("R24").Value = 324.267280999074
$ws.Range("S24").Value = 0.009707360549368538
$ws.Range("T24").Value = 0.01125165264838245

# Row 25
$ws.Range("G25").Value = 2.398207
$ws.Range("H25").Value = 7.194621
$ws.Range("I25").Value = 0.06426264320896187
$ws.Range("J25").Value = 0.06860517894832399
$ws.Range("M25").Value = 23.556204
$ws.Range("N25").Value = 47.112408
$ws.Range("O25").Value = 0.236850294013169
$ws.Range("P25").Value = 0.1714350065796238
$ws.Range("Q25").Value = 56.49265332622799
$ws.Range("R25").Value = 338.955919957368
$ws.Range("S25").Value = 0.015220625938106
$ws.Range("T25").Value = 0.0117613293044022

# Row 26
$ws.Range("G26").Value = 2.398207
$ws.Range("H26").Value = 7.194621
$ws.Range("I26").Value = 0.06426264320896187
$ws.Range("J26").Value = 0.06860517894832399
$ws.Range("M26").Value = 13.62041
$ws.Range("N26").Value = 40.86123000000001
$ws.Range("O26").Value = 0.1369489801107134
$ws.Range("P26").Value = 0.1486879047638899
$ws.Range("Q26").Value = 32.66456260487
$ws.Range("R26").Value = 293.98106344383
$ws.Range("S26").Value = 0.00880070344668599
$ws.Range("T26").Value = 0.01020076031377802
